$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.679.89"
$ws.Range("E2").Value = "  -3.57%  "
$ws.Range("D3").Value = "1.742.19"
$ws.Range("E3").Value = "  -5.53%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'238.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.72%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.4928"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.60%  "
$ws.Range("D8").Value = "'41.55"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.88%  "
$ws.Range("D9").Value = "'0.2429"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -23.17%  "
$ws.Range("D10").Value = "'0.05970"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -12.30%  "
$ws.Range("D11").Value = "1.746.62"
$ws.Range("E11").Value = "  -5.28%  "
$ws.Range("D12").Value = "'0.06806"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -12.64%  "
$ws.Range("E13").Value = "  -22.62%  "
$ws.Range("D14").Value = "'4.458"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -11.16%  "
$ws.Range("D15").Value = "'77.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -12.86%  "
$ws.Range("D16").Value = "'0.5805"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -26.20%  "
$ws.Range("D17").Value = "'0.9990"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "'0.9995"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "25.742.88"
$ws.Range("E19").Value = "  -3.46%  "
$ws.Range("D20").Value = "'11.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -17.30%  "
$ws.Range("D21").Value = "'0.000006455"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -18.59%  "
$ws.Range("D22").Value = "1.963.39"
$ws.Range("E22").Value = "  -5.52%  "
$ws.Range("D23").Value = "'3.978"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -13.82%  "
$ws.Range("D24").Value = "'5.015"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -16.51%  "
$ws.Range("D25").Value = "'7.858"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -15.91%  "
$ws.Range("D26").Value = "'136.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.84%  "
$ws.Range("D27").Value = "'1.479"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -12.11%  "
$ws.Range("D28").Value = "'1.849"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -16.97%  "
$ws.Range("E29").Value = "  -14.84%  "
$ws.Range("D30").Value = "'100.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.40%  "
$ws.Range("D31").Value = "'3.770"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.72%  "
$ws.Range("D32").Value = "'0.08100"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.02%  "
$ws.Range("D33").Value = "'3.359"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -17.86%  "
$ws.Range("D34").Value = "'0.04400"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.83%  "
$ws.Range("D35").Value = "'0.9989"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -8.01%  "
$ws.Range("E37").Value = "  -10.59%  "
$ws.Range("D38").Value = "'0.6106"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -16.47%  "
$ws.Range("D39").Value = "'2.712"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -12.92%  "
$ws.Range("E40").Value = "  -12.28%  "
$ws.Range("D41").Value = "'0.9992"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "'103.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.53%  "
$ws.Range("D43").Value = "'0.01500"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -13.69%  "
$ws.Range("D44").Value = "'0.7750"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -14.65%  "
$ws.Range("D45").Value = "'5.158"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -12.93%  "
$ws.Range("D46").Value = "'0.3758"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -22.37%  "
$ws.Range("D47").Value = "'0.05113"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -12.24%  "
$ws.Range("D50").Value = "'30.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -13.70%  "
$ws.Range("D51").Value = "'52.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -12.22%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1071"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -14.01%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'5.972"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -22.90%  "
